$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of the existing header cell H1 onto the two new header cells
# I1 and J1 so they match the formatting (bold, border, centered) of the
# other header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I (I0) and J (IF) for rows 2-6
$values = @(
    @(1, 5),
    @(1, 5),
    @(6, 7),
    @(1, 3),
    @(1, 2)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
